# "All in Working Condition"
#
# Sheet3 currently only has the header row (Name / Country / Company).
# Sheet1 holds the authoritative Name/Country list (rows 2:61) and
# Sheet2 holds a Name/Company list. Fill in Sheet3 by pulling the
# Name + Country straight from Sheet1 and looking up the matching
# Company from Sheet2, then freeze the lookup results down to plain
# values so the sheet is fully self-contained ("in working condition").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$lastRow = 61

# Bring over Name (A) and Country (B) from Sheet1 into Sheet3.
$namesAndCountries = $ws1.Range("A2:B$lastRow")
$namesAndCountries.Copy()
$ws3.Range("A2").PasteSpecial(-4163)

# Look up each person's Company from Sheet2 (Name in column A, Company in column B).
$sheet2Name = $ws2.Name
for ($r = 2; $r -le $lastRow; $r++) {
    $ws3.Cells.Item($r, 3).Formula = "=VLOOKUP(A$r,$sheet2Name!`$A:`$B,2,FALSE)"
}

# Convert the lookup formulas into static values.
$companyRange = $ws3.Range("C2:C$lastRow")
$companyRange.Copy()
$companyRange.PasteSpecial(-4163)
